$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.411.39"
$ws.Range("E2").Value = "  -3.04%  "

$ws.Range("D3").Value = "1.739.70"
$ws.Range("E3").Value = "  -3.76%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.003"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "321.54"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -4.87%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +0.14%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4238"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -9.49%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3580"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -6.15%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "45.41"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +0.29%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.07410"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -3.03%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "1.109"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -3.85%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +0.10%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "21.34"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -5.09%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "6.076"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -4.05%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "7.173"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -3.73%  "

$ws.Range("D16").Value = "1.740.59"
$ws.Range("E16").Value = "  -3.75%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.00001064"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -2.83%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "86.99"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +6.13%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.06201"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -7.74%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +0.13%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "16.82"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -3.54%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "6.092"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -5.31%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "0.5244"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -5.44%  "

$ws.Range("D24").Value = "27.447.68"
$ws.Range("E24").Value = "  -2.91%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "11.60"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -2.39%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "2.319"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -3.77%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "20.35"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -2.05%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "151.64"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -1.73%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "2.341"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -1.53%  "

$ws.Range("D30").Value = "1.938.03"
$ws.Range("E30").Value = "  -3.74%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "126.24"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -5.11%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "1.206"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -4.15%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "5.660"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -3.44%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.09137"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -5.23%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "3.667"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -9.05%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "12.64"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +4.24%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.02282"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -3.07%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.2127"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -5.93%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "5.069"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -3.73%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.06059"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -5.18%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.6390"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -3.75%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "1.192"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -3.78%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "1.417"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -5.25%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "1.000"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +0.07%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "7.881"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -4.84%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "13.67"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -3.86%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "3.720"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -3.72%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.5866"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -4.51%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "124.97"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -4.44%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "1.950"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -4.35%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.06848"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -4.32%  "
